# Apply updated cryptocurrency price/volume figures (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.088.47"
$ws.Range("E2").Value = "  +2.54%  "
$ws.Range("D3").Value = "1.678.28"
$ws.Range("E3").Value = "  +3.71%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'216.30"
$ws.Range("E5").Value = "  +1.55%  "
$ws.Range("E6").Value = "  +2.02%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  +3.17%  "
$ws.Range("E9").Value = "  +1.85%  "
$ws.Range("D10").Value = "'20.23"
$ws.Range("E10").Value = "  +5.44%  "
$ws.Range("D11").Value = "'0.0888"
$ws.Range("E11").Value = "  +4.83%  "
$ws.Range("D12").Value = "1.913.57"
$ws.Range("E12").Value = "  +3.69%  "
$ws.Range("D13").Value = "1.679.41"
$ws.Range("E13").Value = "  +3.80%  "
$ws.Range("D14").Value = "'4.10"
$ws.Range("E14").Value = "  +1.76%  "
$ws.Range("E15").Value = "  +2.91%  "
$ws.Range("D16").Value = "'65.97"
$ws.Range("E16").Value = "  +3.28%  "
$ws.Range("D17").Value = "27.110.67"
$ws.Range("E17").Value = "  +2.57%  "
$ws.Range("D18").Value = "'238.59"
$ws.Range("E18").Value = "  +0.63%  "
$ws.Range("E19").Value = "  +1.78%  "
$ws.Range("D20").Value = "'7.76"
$ws.Range("E20").Value = "  -0.80%  "
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("E22").Value = "  +4.44%  "
$ws.Range("E23").Value = "  +3.22%  "
$ws.Range("E24").Value = "  +2.37%  "
$ws.Range("D25").Value = "'145.73"
$ws.Range("E25").Value = "  -0.91%  "
$ws.Range("D26").Value = "'7.16"
$ws.Range("E26").Value = "  +1.65%  "
$ws.Range("E27").Value = "  +0.65%  "
$ws.Range("D28").Value = "'16.02"
$ws.Range("E28").Value = "  +3.49%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.16%  "
$ws.Range("E30").Value = "  +0.49%  "
$ws.Range("E31").Value = "  +1.58%  "
$ws.Range("D32").Value = "'3.33"
$ws.Range("E32").Value = "  +2.50%  "
$ws.Range("D33").Value = "1.479.25"
$ws.Range("E33").Value = "  -3.03%  "
$ws.Range("D34").Value = "'3.12"
$ws.Range("E34").Value = "  +4.91%  "
$ws.Range("E35").Value = "  +5.57%  "
$ws.Range("E36").Value = "  -0.45%  "
$ws.Range("E37").Value = "  +1.51%  "
$ws.Range("D38").Value = "'0.903"
$ws.Range("E38").Value = "  +8.73%  "
$ws.Range("E39").Value = "  +2.44%  "
$ws.Range("D40").Value = "'6.08"
$ws.Range("E40").Value = "  +2.31%  "
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("E42").Value = "  +10.60%  "
$ws.Range("D43").Value = "'66.79"
$ws.Range("E43").Value = "  +8.52%  "
$ws.Range("D44").Value = "'2.26"
$ws.Range("E44").Value = "  +3.10%  "
$ws.Range("D45").Value = "1.824.02"
$ws.Range("E45").Value = "  +3.84%  "
$ws.Range("E46").Value = "  +2.20%  "
$ws.Range("D47").Value = "'90.42"
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("E48").Value = "  +2.18%  "
$ws.Range("E49").Value = "  +4.75%  "
$ws.Range("E50").Value = "  +1.16%  "
$ws.Range("D51").Value = "'7.65"
$ws.Range("E51").Value = "  +2.21%  "
